# Generate Report for Handback
# The de0ea274-... file has now been handed back (in sync with en-US).
# Update its Status on every sheet and record the new Handback datetime
# on the per-locale (zh-cn / de-de) report sheets.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: row for de0ea274-...md is row 3 (columns B & C) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Cells.Item(3, 2).Value = $statusText   # B3
$overview.Cells.Item(3, 3).Value = $statusText   # C3

# --- zh-cn sheet: row for de0ea274-...md is row 3 ---
# Column B = Status, Column G = Latest Handback DateTime
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Cells.Item(3, 2).Value = $statusText              # B3 Status
$zhcn.Cells.Item(3, 7).Value = "2016-02-24 08:59:32"    # G3 Latest Handback DateTime

# --- de-de sheet: row for de0ea274-...md is row 3 ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Cells.Item(3, 2).Value = $statusText              # B3 Status
$dede.Cells.Item(3, 7).Value = "2016-02-24 08:59:54"    # G3 Latest Handback DateTime
